# Programming Puzzles workbook - add graph visualiser note + more graph algorithms
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leetcode")

# --- Row 86: Rotting Oranges -> topics refined, mark as STRUGGLED ---
$ws.Cells.Item(86, 4).Value = "Graphs, BFS"        # D86 TOPICS
$ws.Cells.Item(86, 7).Value = "STRUGGLED"          # G86 STATUS

# --- Row 87: Walls and Gates -> topics refined, mark SOLVED w/ date + note ---
$ws.Cells.Item(87, 4).Value = "Graphs, BFS"        # D87 TOPICS
$ws.Cells.Item(87, 7).Value = "SOLVED"             # G87 STATUS
$ws.Cells.Item(87, 8).Value = "25/07/2025"         # H87 LAST SOLVED
$ws.Cells.Item(87, 9).Value = "Similar to Rotting Oranges above. Reverse thinking." # I87 NOTES
$ws.Cells.Item(87, 9).WrapText = $true

# --- Row 88: Course Schedule -> topics refined, mark SOLVED w/ date + note ---
$ws.Cells.Item(88, 4).Value = "Graphs, DFS, BFS"   # D88 TOPICS
$ws.Cells.Item(88, 7).Value = "SOLVED"             # G88 STATUS
$ws.Cells.Item(88, 8).Value = "26/07/2025"         # H88 LAST SOLVED
$ws.Cells.Item(88, 9).Value = "Cycle detection."   # I88 NOTES
$ws.Cells.Item(88, 9).WrapText = $true

# --- Row 89: new entry - Course Schedule II ---
$ws.Cells.Item(89, 2).Value = 210                  # B89 IDENTIFIER
$ws.Cells.Item(89, 2).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(89, 3).Value = "Course Schedule II" # C89 PROBLEM
$ws.Cells.Item(89, 4).Value = "Topological Sort"   # D89 TOPICS
$ws.Cells.Item(89, 5).Value = "Medium"             # E89 DIFFICULTY

# --- Row 90: new entry - Redundant Connection ---
$ws.Cells.Item(90, 2).Value = 684                  # B90 IDENTIFIER
$ws.Cells.Item(90, 2).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(90, 3).Value = "Redundant Connection" # C90 PROBLEM
$ws.Cells.Item(90, 4).Value = "Union Find"         # D90 TOPICS

# --- Row 91: new entry - Number of Connected Components in Undirected Graph ---
$ws.Cells.Item(91, 2).Value = 323                  # B91 IDENTIFIER
$ws.Cells.Item(91, 2).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(91, 3).Value = "Number of Connected Components in Undirected Graph" # C91 PROBLEM
$ws.Cells.Item(91, 4).Value = "Union Find"         # D91 TOPICS

# --- Update active selection to C92 (next empty row, ready for the next entry) ---
$ws.Range("C92").Select() | Out-Null

Write-Output "edit complete"
